# Update countries & provincias Spain
# - Reorder "Camerun" ahead of "Eslovaquia"/"Cuba" in the country list (it
#   overtook them in total cases), shifting the data that used to sit at
#   rows 80/81 down to rows 81/82, and giving Camerun (now row 80) a fresh
#   set of figures.
# - Refresh a few other countries' daily figures (Estados Unidos, Austria).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Estados Unidos (row 4): refresh daily figures
$ws.Range("B4").Value = 673215
$ws.Range("C4").Value = 25067
$ws.Range("E4").Value = 581599
$ws.Range("G4").Value = 1796
$ws.Range("H4").Value = 34384

# Austria (row 20): refresh daily figures
$ws.Range("B20").Value = 14476
$ws.Range("C20").Value = 126
$ws.Range("E20").Value = 5080

# Rows 80-82: Camerun moves above Eslovaquia and Cuba in rank order.
# Row 80 becomes Camerun with new figures; Eslovaquia's and Cuba's previous
# figures shift down one row each.
$ws.Range("A80").Value = "Camerun"
$ws.Range("B80").Value = 996
$ws.Range("C80").Value = 148
$ws.Range("D80").Value = 164
$ws.Range("E80").Value = 810
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 5
$ws.Range("H80").Value = 22

$ws.Range("A81").Value = "Eslovaquia"
$ws.Range("B81").Value = 977
$ws.Range("C81").Value = 114
$ws.Range("D81").Value = 167
$ws.Range("E81").Value = 802
$ws.Range("F81").Value = 5
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 8

$ws.Range("A82").Value = "Cuba"
$ws.Range("B82").Value = 862
$ws.Range("C82").Value = 48
$ws.Range("D82").Value = 171
$ws.Range("E82").Value = 664
$ws.Range("F82").Value = 16
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 27
